{"js": "const replacements = [\n  [\"324\u00d76=\", \"581\u00d73=\"],\n  [\"866\u00d77=\", \"385\u00d72=\"],\n  [\"877\u00d76=\", \"961\u00d79=\"],\n  [\"122\u00d77=\", \"366\u00d78=\"],\n  [\"735\u00d78=\", \"768\u00d79=\"],\n  [\"254\u00d79=\", \"927\u00d76=\"],\n  [\"976\u00d76=\", \"771\u00d74=\"],\n  [\"566\u00d72=\", \"468\u00d76=\"],\n  [\"672\u00d72=\", \"963\u00d72=\"],\n  [\"229\u00d74=\", \"396\u00d76=\"],\n  [\"982\u00d77=\", \"841\u00d78=\"],\n  [\"121\u00d76=\", \"892\u00d76=\"],\n  [\"373\u00d78=\", \"518\u00d79=\"],\n  [\"163\u00d77=\", \"395\u00d75=\"],\n  [\"359\u00d77=\", \"960\u00d78=\"],\n  [\"667\u00d72=\", \"508\u00d77=\"],\n  [\"750\u00d74=\", \"632\u00d79=\"],\n  [\"496\u00d73=\", \"422\u00d77=\"],\n  [\"432\u00d76=\", \"203\u00d75=\"],\n  [\"777\u00d76=\", \"843\u00d73=\"],\n  [\"347\u00d77=\", \"646\u00d73=\"],\n  [\"280\u00d72=\", \"594\u00d72=\"],\n  [\"455\u00d76=\", \"686\u00d73=\"],\n  [\"649\u00d79=\", \"648\u00d72=\"],\n  [\"691\u00d72=\", \"571\u00d72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"324\u00d76=\"\n$find.Replacement.Text = \"581\u00d73=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"866\u00d77=\"\n$find.Replacement.Text = \"385\u00d72=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"877\u00d76=\"\n$find.Replacement.Text = \"961\u00d79=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"122\u00d77=\"\n$find.Replacement.Text = \"366\u00d78=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"735\u00d78=\"\n$find.Replacement.Text = \"768\u00d79=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"254\u00d79=\"\n$find.Replacement.Text = \"927\u00d76=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"976\u00d76=\"\n$find.Replacement.Text = \"771\u00d74=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"566\u00d72=\"\n$find.Replacement.Text = \"468\u00d76=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"672\u00d72=\"\n$find.Replacement.Text = \"963\u00d72=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"229\u00d74=\"\n$find.Replacement.Text = \"396\u00d76=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"982\u00d77=\"\n$find.Replacement.Text = \"841\u00d78=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"121\u00d76=\"\n$find.Replacement.Text = \"892\u00d76=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"373\u00d78=\"\n$find.Replacement.Text = \"518\u00d79=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"163\u00d77=\"\n$find.Replacement.Text = \"395\u00d75=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"359\u00d77=\"\n$find.Replacement.Text = \"960\u00d78=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"667\u00d72=\"\n$find.Replacement.Text = \"508\u00d77=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"750\u00d74=\"\n$find.Replacement.Text = \"632\u00d79=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"496\u00d73=\"\n$find.Replacement.Text = \"422\u00d77=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"432\u00d76=\"\n$find.Replacement.Text = \"203\u00d75=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"777\u00d76=\"\n$find.Replacement.Text = \"843\u00d73=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"347\u00d77=\"\n$find.Replacement.Text = \"646\u00d73=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"280\u00d72=\"\n$find.Replacement.Text = \"594\u00d72=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"455\u00d76=\"\n$find.Replacement.Text = \"686\u00d73=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"649\u00d79=\"\n$find.Replacement.Text = \"648\u00d72=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"691\u00d72=\"\n$find.Replacement.Text = \"571\u00d72=\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
